# Updated cryptos list on Mon Jul 22 18:49:07 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.439.37"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.473.10"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.89"
$ws.Range("E5").Value = "  +0.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.24"
$ws.Range("E6").Value = "  +3.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  +5.83%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.469.54"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10
$ws.Range("E10").Value = "  +9.08%  "

# Row 11
$ws.Range("E11").Value = "  -1.03%  "

# Row 12
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.070.78"
$ws.Range("E13").Value = "  +0.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.01"
$ws.Range("E14").Value = "  +4.68%  "

# Row 15
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.406.24"
$ws.Range("E16").Value = "  +0.65%  "

# Row 17
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.474.15"
$ws.Range("E18").Value = "  +0.38%  "

# Row 19
$ws.Range("E19").Value = "  -0.20%  "

# Row 20
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.26"
$ws.Range("E21").Value = "  +2.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.96"
$ws.Range("E22").Value = "  +1.69%  "

# Row 23
$ws.Range("E23").Value = "  +1.50%  "

# Row 24
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  +1.58%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  -0.84%  "

# Row 27
$ws.Range("E27").Value = "  +1.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.39"
$ws.Range("E28").Value = "  +1.30%  "

# Row 29
$ws.Range("E29").Value = "  -1.22%  "

# Row 30
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("E33").Value = "  +1.30%  "

# Row 34
$ws.Range("E34").Value = "  +0.78%  "

# Row 35
$ws.Range("E35").Value = "  +1.17%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("E37").Value = "  -1.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.02"
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.892"
$ws.Range("E39").Value = "  +2.44%  "

# Row 40
$ws.Range("E40").Value = "  +11.06%  "

# Row 41
$ws.Range("E41").Value = "  -2.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.77"
$ws.Range("E42").Value = "  -1.94%  "

# Row 43
$ws.Range("E43").Value = "  +1.64%  "

# Row 44
$ws.Range("E44").Value = "  +0.38%  "

# Row 45
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.30"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.744.66"
$ws.Range("E47").Value = "  -1.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.54"
$ws.Range("E48").Value = "  -1.28%  "

# Row 49
$ws.Range("E49").Value = "  +0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "326.37"
$ws.Range("E50").Value = "  -3.23%  "

# Row 51
$ws.Range("E51").Value = "  -2.00%  "

